$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

Set-TextValue 'D2' '63.950.17'
Set-TextValue 'E2' '  +4.70%  '
Set-TextValue 'D3' '2.754.58'
Set-TextValue 'E3' '  +3.71%  '
Set-TextValue 'E4' '  +0.28%  '
Set-TextValue 'D5' '580.97'
Set-TextValue 'D6' '155.46'
Set-TextValue 'E6' '  +7.35%  '
Set-TextValue 'D7' '0.997'
Set-TextValue 'E7' '  -0.04%  '
Set-TextValue 'E8' '  +1.65%  '
Set-TextValue 'D9' '2.770.81'
Set-TextValue 'E9' '  +3.67%  '
Set-TextValue 'D10' '6.74'
Set-TextValue 'E10' '  +1.63%  '
Set-TextValue 'D11' '0.113'
Set-TextValue 'E11' '  +4.97%  '
Set-TextValue 'D12' '0.391'
Set-TextValue 'E12' '  +2.20%  '
Set-TextValue 'E13' '  +2.85%  '
Set-TextValue 'D14' '3.249.92'
Set-TextValue 'E14' '  +4.00%  '
Set-TextValue 'D15' '26.78'
Set-TextValue 'E15' '  +2.89%  '
Set-TextValue 'D16' '63.905.87'
Set-TextValue 'E16' '  +4.64%  '
Set-TextValue 'E17' '  +6.32%  '
Set-TextValue 'D18' '2.771.06'
Set-TextValue 'E18' '  +4.02%  '
Set-TextValue 'D19' '12.00'
Set-TextValue 'E19' '  +2.95%  '
Set-TextValue 'D20' '4.90'
Set-TextValue 'E20' '  +3.51%  '
Set-TextValue 'D21' '362.46'
Set-TextValue 'E21' '  +3.22%  '
Set-TextValue 'D22' '7.01'
Set-TextValue 'E22' '  +0.73%  '
Set-TextValue 'E23' '  -0.10%  '
Set-TextValue 'E24' '  +0.60%  '
Set-TextValue 'D25' '66.46'
Set-TextValue 'E25' '  +3.57%  '
Set-TextValue 'D26' '0.172'
Set-TextValue 'E26' '  +5.62%  '
Set-TextValue 'D27' '8.56'
Set-TextValue 'E27' '  +4.54%  '
Set-TextValue 'E28' '  +0.38%  '
Set-TextValue 'D29' '0.0₃0906'
Set-TextValue 'E29' '  +10.95%  '
Set-TextValue 'E30' '  +1.88%  '
Set-TextValue 'D31' '7.14'
Set-TextValue 'E31' '  +3.45%  '
Set-TextValue 'D32' '1.30'
Set-TextValue 'E32' '  +21.71%  '
Set-TextValue 'D33' '172.81'
Set-TextValue 'E33' '  +3.98%  '
Set-TextValue 'D34' '0.998'
Set-TextValue 'E34' '  -0.03%  '
Set-TextValue 'D35' '20.55'
Set-TextValue 'E35' '  +2.97%  '
Set-TextValue 'D36' '4.86'
Set-TextValue 'E36' '  +7.62%  '
Set-TextValue 'E37' '  +8.22%  '
Set-TextValue 'D38' '1.83'
Set-TextValue 'E38' '  +9.22%  '
Set-TextValue 'D39' '1.02'
Set-TextValue 'E39' '  +13.59%  '
Set-TextValue 'D40' '346.06'
Set-TextValue 'E40' '  +1.96%  '
Set-TextValue 'D41' '4.25'
Set-TextValue 'E41' '  +5.03%  '
Set-TextValue 'D42' '39.36'
Set-TextValue 'E42' '  +2.14%  '
Set-TextValue 'D43' '5.89'
Set-TextValue 'E43' '  +11.98%  '
Set-TextValue 'D44' '22.07'
Set-TextValue 'E44' '  +8.13%  '
Set-TextValue 'D45' '22.09'
Set-TextValue 'E45' '  +7.41%  '
Set-TextValue 'D46' '0.0594'
Set-TextValue 'E46' '  +5.79%  '
Set-TextValue 'D47' '0.650'
Set-TextValue 'E47' '  +5.64%  '
Set-TextValue 'D48' '137.77'
Set-TextValue 'E48' '  +2.65%  '
Set-TextValue 'E49' '  +3.11%  '
Set-TextValue 'E50' '  +1.57%  '
Set-TextValue 'D51' '0.996'
Set-TextValue 'E51' '  -0.15%  '
